$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new row 4 values (keywords/amounts)
$ws.Range("B4").Value = -4000
$ws.Range("C4").Value = 8000
$ws.Range("D4").Value = -2000
$ws.Range("E4").Value = -2000

# Recalculate formulas so dependent cells (I2:I5) reflect the new totals
$excel.Calculate()

# Update the active selection to match the authored state
$ws.Range("B5").Select()
